# Add 2022-Q4 data.
#
# The workbook currently has a "2022-Q3" detail sheet (holding fund
# position data) and a "总计" summary sheet with one data row pointing at
# that quarter. We need to:
#   1. Introduce a new "2022-Q4" detail sheet (placed right before the
#      existing "2022-Q3" sheet) carrying the new quarter's fund figures.
#   2. Leave the existing "2022-Q3" sheet and its data exactly as-is.
#   3. Update "总计": its existing single data row now describes 2022-Q4
#      totals, and a new row is appended below it for the 2022-Q3 totals
#      that used to live there.

$wb = $excel.ActiveWorkbook

# --- 1 & 2: create the "2022-Q4" sheet ------------------------------------
# Duplicate "2022-Q3" (copy placed immediately before the original) so the
# new sheet starts out with identical layout/formatting, then rename the
# copy and overwrite its figures with the Q4 numbers. The original
# "2022-Q3" sheet is left untouched.
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3, $null)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# These columns hold plain text in the source data (fund codes, ratios
# formatted as strings, etc.) - force text formatting so values such as
# "166109" or "0.90" are not reinterpreted as numbers.
$q4.Range("B2:G3").NumberFormat = "@"

$q4.Range("B2").Value = "166109"
$q4.Range("C2").Value = "信澳量化先锋混合（LOF）A"
$q4.Range("D2").Value = "0.90"
$q4.Range("E2").Value = "94.26"
$q4.Range("F2").Value = "2.82"
$q4.Range("G2").Value = "0.0254"
$q4.Range("H2").Value = 8

$q4.Range("B3").Value = "166110"
$q4.Range("C3").Value = "信澳量化先锋混合（LOF）C"
$q4.Range("D3").Value = "0.20"
$q4.Range("E3").Value = "94.26"
$q4.Range("F3").Value = "2.82"
$q4.Range("G3").Value = "0.0056"
$q4.Range("H3").Value = 8

# --- 3: update the "总计" summary sheet -----------------------------------
$total = $wb.Worksheets.Item("总计")

# Seed row 3 with row 2's formatting (so the new A3 cell matches A2's
# style) before overwriting the values.
$total.Range("A2").Copy($total.Range("A3"))

$total.Range("B2").Value = "2022-Q4"
$total.Range("D2").Value = 0.03

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.02

# Keep the originally-selected "2022-Q3" tab active, same as before the edit
# (the brand-new "2022-Q4" sheet is not the selected tab). Re-fetch the
# worksheet reference since it shifted position during the copy above.
$wb.Worksheets.Item("2022-Q3").Activate()
